$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 10
$ws.Range("L2").Value = 1.33
$ws.Range("M2").Value = 3.4
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 1.8

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 2.9
$ws.Range("K3").Value = 9
$ws.Range("N3").Value = 2.2
$ws.Range("O3").Value = 1.67
$ws.Range("AE3").Value = 9.5

# Row 4
$ws.Range("H4").Value = 3.4
$ws.Range("T4").Value = 9

# Row 6
$ws.Range("G6").Value = 2.05
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 3.5
$ws.Range("L6").Value = 1.22
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 1.73
$ws.Range("O6").Value = 2.08
$ws.Range("V6").Value = 9
$ws.Range("Y6").Value = 23
$ws.Range("AB6").Value = 13
$ws.Range("AE6").Value = 12
$ws.Range("AF6").Value = 19
$ws.Range("AG6").Value = 12

# Row 10
$ws.Range("K10").Value = 13

# Row 11
$ws.Range("P11").Value = 1.36
$ws.Range("Q11").Value = 3
$ws.Range("W11").Value = 13
$ws.Range("Z11").Value = 12
$ws.Range("AC11").Value = 51
$ws.Range("AJ11").Value = 41

# Row 15
$ws.Range("G15").Value = 2.38
$ws.Range("I15").Value = 2.9
$ws.Range("AH15").Value = 29

# Row 19
$ws.Range("K19").Value = 9

# Row 24
$ws.Range("L24").Value = 1.44
$ws.Range("M24").Value = 2.63

# Row 25
$ws.Range("N25").Value = 1.98
$ws.Range("O25").Value = 1.88
